# Update countries & provincias Spain
#
# - Refreshed case counts for Estados Unidos, India, Australia and Crucero.
# - Inserted two "new" countries into the ranked (by Casos totales, desc)
#   list: Hungria (ahead of Islandia/Barein) and Bulgaria (ahead of Tunez).
#   Because the sheet is sorted by total cases, adding these rows pushes the
#   countries that used to occupy that rank down by one row, which is
#   expressed here as a 3-row rotation of values rather than a physical
#   row insert (no shift of the rows below is required/observed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 678210
$ws.Range("C4").Value = 640
$ws.Range("E4").Value = 585725

# --- India (row 21) ---
$ws.Range("B21").Value = 13495
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 1777
$ws.Range("E21").Value = 11270

# --- Australia (row 35) ---
$ws.Range("B35").Value = 6497
$ws.Range("C35").Value = 29
$ws.Range("E35").Value = 2687

# --- Hungria / Islandia / Barein block (rows 62-64) ---
# Hungria is newly promoted into row 62 with fresh data; the former row 62
# (Islandia) and row 63 (Barein) data moves down one row each.
$ws.Range("A62").Value = "Hungria"
$ws.Range("B62").Value = 1763
$ws.Range("C62").Value = 111
$ws.Range("D62").Value = 207
$ws.Range("E62").Value = 1400
$ws.Range("F62").Value = 63
$ws.Range("G62").Value = 14
$ws.Range("H62").Value = 156

$ws.Range("A63").Value = "Islandia"
$ws.Range("B63").Value = 1739
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 1144
$ws.Range("E63").Value = 587
$ws.Range("F63").Value = 6
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 8

$ws.Range("A64").Value = "Barein"
$ws.Range("B64").Value = 1700
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 703
$ws.Range("E64").Value = 990
$ws.Range("F64").Value = 3
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 7

# --- Bulgaria / Tunez block (rows 84-85) ---
# Bulgaria is newly promoted into row 84 with fresh data; the former row 84
# (Tunez) data moves down one row.
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 825
$ws.Range("C84").Value = 25
$ws.Range("D84").Value = 141
$ws.Range("E84").Value = 644
$ws.Range("F84").Value = 37
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 40

$ws.Range("A85").Value = "Tunez"
$ws.Range("B85").Value = 822
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 43
$ws.Range("E85").Value = 742
$ws.Range("F85").Value = 89
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 37

# --- Crucero (row 87) ---
$ws.Range("E87").Value = 55
$ws.Range("H87").Value = 13
